$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $orig = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $orig
}

Set-TextValue 'D2' '60.698.05'
Set-TextValue 'E2' '  -1.13%  '
Set-TextValue 'D3' '3.365.42'
Set-TextValue 'E3' '  +1.15%  '
Set-TextValue 'E4' '  +0.29%  '
Set-TextValue 'D5' '403.39'
Set-TextValue 'E5' '  -1.58%  '
Set-TextValue 'D6' '127.32'
Set-TextValue 'E6' '  +11.88%  '
Set-TextValue 'D7' '0.601'
Set-TextValue 'E7' '  +5.82%  '
Set-TextValue 'E8' '  +0.14%  '
Set-TextValue 'E9' '  +7.45%  '
Set-TextValue 'E10' '  +12.29%  '
Set-TextValue 'D11' '41.79'
Set-TextValue 'E11' '  +7.44%  '
Set-TextValue 'E12' '  -0.90%  '
Set-TextValue 'D13' '3.909.54'
Set-TextValue 'E13' '  +3.44%  '
Set-TextValue 'D14' '8.45'
Set-TextValue 'E14' '  +3.27%  '
Set-TextValue 'D15' '19.49'
Set-TextValue 'E15' '  +2.60%  '
Set-TextValue 'D16' '3.364.86'
Set-TextValue 'E16' '  +0.80%  '
Set-TextValue 'B17' 'WrappedBTC'
Set-TextValue 'C17' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D17' '60.800.04'
Set-TextValue 'E17' '  -0.38%  '
Set-TextValue 'B18' 'Uniswap'
Set-TextValue 'C18' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D18' '11.33'
Set-TextValue 'E18' '  +8.20%  '
Set-TextValue 'E19' '  +2.29%  '
Set-TextValue 'E20' '  +18.16%  '
Set-TextValue 'D21' '3.23'
Set-TextValue 'E21' '  +0.51%  '
Set-TextValue 'D22' '82.56'
Set-TextValue 'E22' '  +12.51%  '
Set-TextValue 'D23' '13.08'
Set-TextValue 'E23' '  +5.49%  '
Set-TextValue 'D24' '304.57'
Set-TextValue 'E24' '  +3.07%  '
Set-TextValue 'E25' '  +2.22%  '
Set-TextValue 'B26' 'LEO'
Set-TextValue 'C26' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D26' '4.72'
Set-TextValue 'E26' '  +5.30%  '
Set-TextValue 'B27' 'Filecoin'
Set-TextValue 'C27' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D27' '8.41'
Set-TextValue 'E27' '  +12.91%  '
Set-TextValue 'D28' '29.39'
Set-TextValue 'E28' '  +2.36%  '
Set-TextValue 'D29' '7.42'
Set-TextValue 'E29' '  +0.27%  '
Set-TextValue 'E30' '  +1.84%  '
Set-TextValue 'E31' '  +5.51%  '
Set-TextValue 'D32' '11.75'
Set-TextValue 'E32' '  +5.34%  '
Set-TextValue 'E33' '  +6.87%  '
Set-TextValue 'D34' '42.10'
Set-TextValue 'E34' '  +4.77%  '
Set-TextValue 'E35' '  +0.08%  '
Set-TextValue 'D36' '0.0482'
Set-TextValue 'E36' '  +1.04%  '
Set-TextValue 'D37' '52.21'
Set-TextValue 'D38' '0.996'
Set-TextValue 'E38' '  +0.06%  '
Set-TextValue 'D39' '3.40'
Set-TextValue 'E39' '  +2.81%  '
Set-TextValue 'D40' '2.96'
Set-TextValue 'E40' '  -3.08%  '
Set-TextValue 'D41' '2.03'
Set-TextValue 'E41' '  +7.84%  '
Set-TextValue 'E42' '  +4.42%  '
Set-TextValue 'D43' '135.87'
Set-TextValue 'E43' '  +0.38%  '
Set-TextValue 'D44' '3.91'
Set-TextValue 'E44' '  +3.86%  '
Set-TextValue 'D45' '16.84'
Set-TextValue 'E45' '  +4.19%  '
Set-TextValue 'D46' '0.281'
Set-TextValue 'E46' '  -0.29%  '
Set-TextValue 'E47' '  +1.21%  '
Set-TextValue 'D48' '21.72'
Set-TextValue 'E48' '  +3.70%  '
Set-TextValue 'D49' '2.131.87'
Set-TextValue 'E49' '  +0.81%  '
Set-TextValue 'D50' '3.706.53'
Set-TextValue 'E50' '  +1.79%  '
Set-TextValue 'D51' '2.35'
Set-TextValue 'E51' '  +0.98%  '
